# Fashion Image Tagging - content update
$d = $word.ActiveDocument
$report = @()

# (1) Rewrite the opening "Inspiration" paragraph.
$ok = $d.Content.Find.Execute("In this digital age, people are constantly seeking ways to simplify their lives leveraging technology. This is where closet apps come in. ", $true, $false, $false, $false, $false, $true, 1, $false, "Time is precious and we do not always have the time to look nice. Sometimes, you just throw on whatever you can find in your closet without a second thought. That is where closet apps come in. There are many closets apps available on mobile and digital platform that offer us the promise to organize your wardrobe for easy visibility, plan daily outfits, and track clothing use to declutter your closet by identifying clothes rarely worn. ", 2)
$report += "intro-rewrite=$ok"

# (2) Rewrite the follow-up paragraph about closet apps / Fashion-MNIST teaser.
$ok = $d.Content.Find.Execute("There are 100s of closets apps available on major mobile and digital platforms, which offer you the promise of organizing your clothes for easy visibility, plan outfits for months on end, track clothing use to help calculate your cost/wear, and declutter your closet by identifying outfits rarely worn for easy disposal. An examination of the top 5 closet apps in the Apple App store offer the ability to mass upload pictures of clothing and deliver on the aforementioned features; however, there is a fundamental opportunity yet to be explored. In order to get the app ready, the user has to manually tag and classify each item and separate them into respective folders. What if we can use Machine Learning to automate this process and get the application ready for use in minutes versus hours? Wouldn’t this be a competitive advantage? Cue in Fashion-MNIST!", $true, $false, $false, $false, $false, $true, 1, $false, "However, when examining the top 5 closet apps in the Apple App store, the ability to mass upload pictures of clothing has yet to be explored. In order to initialize the app, the user has to manually classify each item via tags then separate them into their respective clothing category. But what if we can use Machine Learning to automate this process and prepare the application in minutes rather than hours? Cue in Fashion-MNIST!", 2)
$report += "closet-apps-rewrite=$ok"

# (3) No textual change, but collapse the three runs that make up this sentence into one
#     (re-assert the same text through Find/Replace so Word re-flows it as a single run).
$ok = $d.Content.Find.Execute("The opportunity to explore machine learning opportunities based on fashion classification has been one of great interest to the scientific community. For this reason, two prominent professors in the data science & analytics space leveraged the idea of the well-loved and vastly used handwritten digit dataset, MNIST, and developed their a new MNIST specifically for Fashion clothing identification. ", $true, $false, $false, $false, $false, $true, 1, $false, "The opportunity to explore machine learning opportunities based on fashion classification has been one of great interest to the scientific community. For this reason, two prominent professors in the data science & analytics space leveraged the idea of the well-loved and vastly used handwritten digit dataset, MNIST, and developed their a new MNIST specifically for Fashion clothing identification. ", 2)
$report += "mnist-origin-merge=$ok"

# (4) Same idea for the italic pull-quote (merge the three runs into one run).
$ok = $d.Content.Find.Execute("“Fashion-MNIST is a dataset of Zalando’s article images consisting of a training set of 60,000 examples and a test set of 10,000 examples. Each example is a 28×28 grayscale image, associated with a label from 10 classes. Fashion-MNIST is intended to serve as a direct drop-in replacement of the original MNIST dataset for benchmarking machine learning algorithms”", $true, $false, $false, $false, $false, $true, 1, $false, "“Fashion-MNIST is a dataset of Zalando’s article images consisting of a training set of 60,000 examples and a test set of 10,000 examples. Each example is a 28×28 grayscale image, associated with a label from 10 classes. Fashion-MNIST is intended to serve as a direct drop-in replacement of the original MNIST dataset for benchmarking machine learning algorithms”", 2)
$report += "quote-merge=$ok"

# (5) Mark the inline picture's run as NoProof (adds <w:noProof/> to its run properties).
$picture = $d.InlineShapes.Item(1)
$picture.Range.NoProofing = $true
$report += "noproof=$($picture.Range.NoProofing)"

# (6) Tidy the "real-world application" phrasing.
$ok = $d.Content.Find.Execute("in a real-world application.", $true, $false, $false, $false, $false, $true, 1, $false, "for real word application.", 2)
$report += "real-word-fix=$ok"

# (7) Remove the trailing "Prediction images gallery" caption paragraph entirely.
$deleted = $false
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "Prediction images gallery*") {
        $para.Range.Delete()
        $deleted = $true
    }
}
$report += "gallery-caption-removed=$deleted"

Write-Output ($report -join "; ")
